# Update the dSF column (F) values on Sheet1 with the re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 2
    3  = -2
    4  = 3
    5  = -1
    7  = 3
    9  = -4
    10 = -1
    11 = -3
    13 = 5
    14 = -2
    16 = -2
    17 = -1
    20 = 1
    21 = 2
    22 = 1
    24 = 2
    25 = 2
    26 = -3
    28 = 3
    29 = 1
    30 = -1
    31 = 1
    32 = 5
    35 = 4
    36 = -4
    37 = 3
    38 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
